$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: was 1x6000, now 2x4090
$ws.Range("A6").Value = "2x4090"
$ws.Range("B6").Value = 467.26
$ws.Range("C6").Value = 0.78
$ws.Range("D6").Value = 0.4636961577422991

# Row 7: 2x4090 label unchanged, throughput/price updated
$ws.Range("B7").Value = 4556.7
$ws.Range("D7").Value = 0.04754903036554232

# Row 8: 2x5090 label unchanged, throughput/price updated
$ws.Range("B8").Value = 1230.14
$ws.Range("D8").Value = 0.2935528566757533

# Row 9: was 4x4090, now 2x5090
$ws.Range("A9").Value = "2x5090"
$ws.Range("B9").Value = 8411.16
$ws.Range("C9").Value = 1.3
$ws.Range("D9").Value = 0.04293237925697658

# Row 10: was 4x5090, now 4x4090
$ws.Range("A10").Value = "4x4090"
$ws.Range("B10").Value = 8902.5
$ws.Range("C10").Value = 1.56
$ws.Range("D10").Value = 0.04867546569315735

# Row 11: new row, 4x4090
$ws.Range("A11").Value = "4x4090"
$ws.Range("B11").Value = 906.1900000000001
$ws.Range("C11").Value = 1.56
$ws.Range("D11").Value = 0.4781925791868519

# Row 12: new row, 4x4090
$ws.Range("A12").Value = "4x4090"
$ws.Range("B12").Value = 1731.44
$ws.Range("C12").Value = 1.56
$ws.Range("D12").Value = 0.2502733755332748

# Row 13: was row 11 4x5090, shifted down with updated values
$ws.Range("A13").Value = "4x5090"
$ws.Range("B13").Value = 2501.38
$ws.Range("C13").Value = 2.6
$ws.Range("D13").Value = 0.2887295101992589
